$d = $word.ActiveDocument

# =====================================================================
# Edit 1: on the "09/11/2021" line, a "." was typed in the middle of the
# run of trailing spaces (2 spaces kept before it, 9 left after it), and
# Word's "_GoBack" (last-edit-position) bookmark now sits right after
# the new ".".
# =====================================================================

# Locate "09/11/2021" in the document body.
$dateFind = $d.Content
$null = $dateFind.Find.Execute("09/11/2021")
$dateStart = $dateFind.Start

# The insertion point is right after "09/11/2021" plus the first two of
# the eleven trailing spaces.
$insertPos = $dateStart + 10 + 2

$insPoint = $d.Range($insertPos, $insertPos)
$insPoint.InsertAfter(".")

# Force the newly typed "." to live in its own run (instead of silently
# re-merging with its neighbours) by round-tripping its font colour.
$dotRange = $d.Range($insertPos, $insertPos + 1)
$dotRange.Font.Color = 41087

# Word always keeps a single hidden "_GoBack" bookmark marking the most
# recent edit location; re-adding it here moves it off of its old spot
# (at "condition." near the end of the document) to right after the ".".
$bmPoint = $d.Range($insertPos + 1, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# Restore the "." run's real formatting (same colour as its neighbours).
$dotRange2 = $d.Range($insertPos, $insertPos + 1)
$dotRange2.Font.Color = 15773696

# =====================================================================
# Edit 2: near the end of the document, " les " and "condition." (which
# were split across a grammar-check proofErr range, with the old
# "_GoBack" bookmark sitting right after "condition.") become a single
# run " les condition.". Moving "_GoBack" away above already dropped the
# stray bookmark here; we still need to drop the gramStart/gramEnd
# proofErr markers and merge the two runs.
# =====================================================================

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Find the " les condition." span to merge (probe, non-destructive).
$probe = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$null = $probe.Find.Execute(" les condition.")
$matchStart = $probe.Start

# The trailing gramEnd proofErr marker sits exactly at the end of the
# paragraph's text (right before the paragraph mark), so a replace whose
# range ends there leaves it stranded. Temporarily extend the paragraph
# with a throwaway marker so the replacement span truly crosses over
# that boundary, which drops the now-interior proofErr/bookmark nodes.
$paraEnd = $lastPara.Range.End
$endPoint = $d.Range($paraEnd, $paraEnd)
$endPoint.InsertAfter("ZZZMARK")

$full = $d.Range($matchStart, $lastPara.Range.End)
$null = $full.Find.Execute(" les condition.ZZZMARK", $false, $false, $false, $false, $false, `
                            $true, 1, $false, " les condition.ZZZMARK", 2)

# Remove the throwaway marker again.
$pEnd2 = $lastPara.Range.End
$dummyRange = $d.Range($pEnd2 - 1 - 7, $pEnd2 - 1)
$dummyRange.Delete()
